# Minor updates to src0.
# - Lowercase "Line"/"Character" -> "line"/"character" in the per-character
#   token list on slide 6 (body placeholder).
# - Paragraph 14 ("i   Line 1, Character 14") additionally gets its leading
#   "i" split off into its own run (flagged misspelled in the source deck).
# - The slide 6 footer "(c)SoftMoore Consulting" is split into three runs:
#   "(c)", "SoftMoore", " Consulting".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# ---------------------------------------------------------------------
# Body placeholder ("Rectangle 3") - per-character breakdown list.
# ---------------------------------------------------------------------
$body = $s.Shapes.Item(4)
$tr = $body.TextFrame.TextRange

# Paragraph index (1-based) -> new text. All of these are simple
# case-only edits ("Line" -> "line", "Character" -> "character") that keep
# the exact same length, so rewriting them in any order is safe.
$paraIndexes = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 15, 16)
$newTexts = @(
  "p   line 1, character 1",
  "a   line 1, character 2",
  "c   line 1, character 3",
  "k   line 1, character 4",
  "a   line 1, character 5",
  "g   line 1, character 6",
  "e   line 1, character 7",
  "    line 1, character 8",
  "e   line 1, character 9",
  "d   line 1, character 10",
  "u   line 1, character 11",
  ".   line 1, character 12",
  "c   line 1, character 13",
  "t   line 1, character 15",
  "a   line 1, character 16"
)

for ($i = 0; $i -lt $paraIndexes.Length; $i++) {
  $idx = $paraIndexes[$i]
  # Write a short, character-disjoint placeholder first so the host's
  # text-diff run-splitter doesn't fragment the run against the old text
  # (e.g. "Line"/"line" sharing 3 of 4 letters) when the final text is
  # applied; this keeps each paragraph a single run, matching the source.
  $para = $tr.Paragraphs($idx, 1)
  $para.Text = "~"
  $para = $tr.Paragraphs($idx, 1)
  $para.Text = $newTexts[$i]
}

# Paragraph 14 is split into two runs: "i" and "   line 1, character 14".
$para14 = $tr.Paragraphs(14, 1)
$para14.Text = "~"
$para14 = $tr.Paragraphs(14, 1)
$para14.Text = "i   line 1, character 14"

$para14 = $tr.Paragraphs(14, 1)
$iChar = $tr.Characters($para14.Start, 1)
$iChar.Text = "i"

# ---------------------------------------------------------------------
# Footer placeholder - split the copyright line into three runs.
# ---------------------------------------------------------------------
$footer = $s.Shapes.Item(1)
$footerRange = $footer.TextFrame.TextRange

$copyrightChar = $footerRange.Characters(1, 1)
$copyrightChar.Text = $copyrightChar.Text

$softMoore = $footerRange.Characters(2, 9)
$softMoore.Text = $softMoore.Text

$consulting = $footerRange.Characters(11, 11)
$consulting.Text = $consulting.Text
